# ADD results from server
# Refresh a handful of investment-cost values (row 2 of each year sheet)
# with updated figures produced by the latest server run.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("N2").Value = 70831.95557958097
$ws1.Range("O2").Value = 69610.44223910036

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("B2").Value = 47338.61932520662
$ws2.Range("I2").Value = 283167.7766510943
$ws2.Range("L2").Value = 178095.3756971828
$ws2.Range("M2").Value = 114008.3253427963
$ws2.Range("N2").Value = 33931.82461160053
$ws2.Range("O2").Value = 50657.26889981552

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 29273.60317916481
$ws3.Range("B2").Value = 22330.72247668595
$ws3.Range("E2").Value = 110739.3594843864
$ws3.Range("I2").Value = 150386.9441391908
$ws3.Range("M2").Value = 35556.98862372932
$ws3.Range("N2").Value = 44813.41193308897
$ws3.Range("O2").Value = 26775.55841092002

$ws4 = $wb.Worksheets.Item("2040")
$ws4.Range("N2").Value = 1041.156112142704

$ws5 = $wb.Worksheets.Item("2045")
$ws5.Range("A2").Value = 34256.25291363284
$ws5.Range("N2").Value = 5271.89502409355
$ws5.Range("O2").Value = 22972.54525065989
